$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('L3').Value = '*maa://22880 (66.67), maa://20276 (84.71), *maa://22749 (66.67)'
$ws.Range('X3').Value = 'maa://27396 (84.69), maa://27484 (96.08), maa://27480 (82.86)'
$ws.Range('AB5').Value = '*maa://29863 (70.97), ***maa://22752 (13.33), **maa://26013 (42.86)'
$ws.Range('D6').Value = 'maa://42407 (94.12)'
$ws.Range('A8').Value = '更新日期：2024.12.15 13:19:36'
$ws.Range('P8').Value = 'maa://32931 (85.29), *maa://21916 (61.29), maa://23252 (92.42), maa://37496 (96.3), **maa://22759 (45.45)'
$ws.Range('X8').Value = 'maa://21411 (95.87)'
$ws.Range('P10').Value = 'maa://28977 (91.36), maa://36669 (88.89), *maa://23264 (61.82)'
$ws.Range('T11').Value = 'maa://22747 (93.38), maa://22501 (98.55)'
$ws.Range('H12').Value = 'maa://21867 (90.18)'
$ws.Range('X12').Value = 'maa://22753 (91.3), *maa://21485 (76.87), maa://37962 (88.89)'
$ws.Range('D13').Value = 'maa://24999 (91.74), maa://36673 (92.65), maa://25001 (85.51)'
$ws.Range('D14').Value = 'maa://30764 (88.46)'
$ws.Range('AF15').Value = 'maa://21364 (80.66), *maa://22766 (70.64), *maa://36666 (78.31)'
$ws.Range('H17').Value = 'maa://22430 (88.59), maa://39599 (86.11)'
$ws.Range('AF19').Value = '*maa://21663 (61.54)'
$ws.Range('D20').Value = 'maa://21432 (89.8), maa://25198 (93.0), *maa://20795 (51.18), maa://36680 (96.55)'
$ws.Range('H20').Value = 'maa://22864 (89.04)'
$ws.Range('P20').Value = 'maa://37442 (94.59)'
$ws.Range('D23').Value = '***maa://28036 (27.94), *maa://41753 (58.33)'
$ws.Range('D24').Value = '*maa://24368 (79.6)'
$ws.Range('X24').Value = 'maa://29988 (86.61), maa://23504 (93.23), **maa://22892 (39.58), *maa://25141 (77.6), maa://36663 (81.82), ***maa://22815 (23.08)'
$ws.Range('D25').Value = 'maa://29753 (95.22)'
$ws.Range('L27').Value = 'maa://28071 (88.89)'
$ws.Range('D28').Value = 'maa://24465 (90.77), maa://25725 (83.53)'
$ws.Range('X28').Value = 'maa://39929 (89.58), ***maa://39723 (14.29), maa://41749 (90.0)'
$ws.Range('H32').Value = 'maa://21895 (97.24), maa://36667 (98.33), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (95.77), maa://41108 (87.76), maa://41238 (96.1)'
$ws.Range('P37').Value = 'maa://21280 (88.89), *maa://21239 (72.73)'
$ws.Range('H47').Value = 'maa://27410 (96.17), maa://29661 (97.81), maa://28038 (84.62)'
